$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.058.83'
$ws.Range("E2").Value = '  +3.60%  '
$ws.Range("D3").Value = '1.724.57'
$ws.Range("E3").Value = '  +2.58%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.87'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.08'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +13.03%  '
$ws.Range("E9").Value = '  +3.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0632'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").Value = '1.969.39'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("D13").Value = '1.729.25'
$ws.Range("E13").Value = '  +2.84%  '
$ws.Range("E14").Value = '  +3.24%  '
$ws.Range("E15").Value = '  +4.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.46'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.01%  '
$ws.Range("D17").Value = '28.025.96'
$ws.Range("E17").Value = '  +3.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.20'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("D19").Value = '0.0₃0755'
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.90'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("E22").Value = '  +2.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.67'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.92%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.81'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.50'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.69'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("E28").Value = '  +1.28%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("D33").Value = '1.495.03'
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.65'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.952'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.606'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("E38").Value = '  +1.28%  '
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.52'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.82'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.29%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("D45").Value = '1.873.54'
$ws.Range("E45").Value = '  +2.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.796'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("E47").Value = '  +12.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '90.95'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").Value = '0.0₆0111'
$ws.Range("E49").Value = '  +4.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.21'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.01%  '
$ws.Range("E51").Value = '  +0.33%  '
